$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts N:P -> O:Q),
# matching the "Late" column being pushed right to make room for a
# new (currently unlabeled) column used for Variable Instalments.
$ws.Columns("N:N").Insert()

# The newly inserted column inherits formatting from its neighbour;
# give it the same width as column M ("In Advance" / Due column).
$ws.Columns("N:N").ColumnWidth = 9.83

# Make "Repayment schedule" the active/selected sheet and select Q7,
# where the workbook was left when saved.
$ws.Activate() | Out-Null
$ws.Range("Q7").Select() | Out-Null
